# Actualización automática de tasas-transfi.xlsx

$wb = $excel.ActiveWorkbook

# --- Update Hoja1!A1 text with new conversion rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.92 = 15080.0 pesos`n✅ 15080.0 pesos = 3.9 = 949.16 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Update tasas sheet numeric values ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 255
$ws2.Range("O10").Value = 3845.4
$ws2.Range("N12").Value = 3864.99
$ws2.Range("O12").Value = 243.27
